$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two oldest rows (2008年, 2009年) - this shifts all data rows up by 2
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# After deletion, former row 14 (2020年) is now row 12.
# Append the new 2021年 row as row 13, seeded from row 12's layout
# (so the blank E/F cells and A-column style carry over correctly).
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 24731.8
$ws.Range("C13").Value = 8103.1
$ws.Range("D13").Value = 24642.1
$ws.Range("G13").Value = 4.1

$wb.Save()
